$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    3 = @{ A = 20170926; E = 18 }
    4 = @{ A = 20170927; E = 13 }
    5 = @{ A = 20170928; E = 7 }
    6 = @{ A = 20170929; E = 16 }
    7 = @{ A = 20170930; E = 9 }
    8 = @{ A = 20170931; E = 17 }
    9 = @{ A = 20170932; E = 5 }
    10 = @{ A = 20170933; E = 20 }
    11 = @{ A = 20170934; E = 16 }
    12 = @{ A = 20170935; E = 7 }
    13 = @{ A = 20170936; E = 6 }
    14 = @{ A = 20170937; E = 6 }
    15 = @{ A = 20170938; E = 8 }
    16 = @{ A = 20170939; E = 8 }
    17 = @{ A = 20170940; E = 17 }
    18 = @{ A = 20170941; E = 19 }
    19 = @{ A = 20170942; E = 13 }
    20 = @{ A = 20170943; E = 10 }
    21 = @{ A = 20170944; E = 9 }
    22 = @{ A = 20170945; E = 12 }
    23 = @{ A = 20170946; E = 11 }
    24 = @{ A = 20170947; E = 5 }
    25 = @{ A = 20170948; E = 8 }
    26 = @{ A = 20170949; E = 19 }
    27 = @{ A = 20170950; E = 9 }
    28 = @{ A = 20170951; E = 13 }
    29 = @{ A = 20170952; E = 18 }
    30 = @{ A = 20170953; E = 7 }
    31 = @{ A = 20170954; E = 13 }
    32 = @{ A = 20170955; E = 17 }
    33 = @{ A = 20170956; E = 16 }
    34 = @{ A = 20170957; E = 6 }
    35 = @{ A = 20170958; E = 19 }
    36 = @{ A = 20170959; E = 6 }
    37 = @{ A = 20170960; E = 13 }
    38 = @{ A = 20170961; E = 7 }
    39 = @{ A = 20170962; E = 16 }
    40 = @{ A = 20170963; E = 17 }
    41 = @{ A = 20170964; E = 11 }
    42 = @{ A = 20170965; E = 15 }
    43 = @{ A = 20170966; E = 6 }
    44 = @{ A = 20170967; E = 7 }
    45 = @{ A = 20170968; E = 10 }
    46 = @{ A = 20170969; E = 18 }
    47 = @{ A = 20170970; E = 11 }
    48 = @{ A = 20170971; E = 5 }
    49 = @{ A = 20170972; E = 17 }
    50 = @{ A = 20170973; E = 8 }
    51 = @{ A = 20170974; E = 8 }
    52 = @{ A = 20170975; E = 20 }
    53 = @{ A = 20170976; E = 16 }
    54 = @{ A = 20170977; E = 5 }
    55 = @{ A = 20170978; E = 7 }
    56 = @{ A = 20170979; E = 17 }
    57 = @{ A = 20170980; E = 5 }
    58 = @{ A = 20170981; E = 13 }
    59 = @{ A = 20170982; E = 16 }
    60 = @{ A = 20170983; E = 6 }
    61 = @{ A = 20170984; E = 15 }
    62 = @{ A = 20170985; E = 5 }
    63 = @{ A = 20170986; E = 9 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item([int]$row, 1).Value2 = $vals.A
    $ws.Cells.Item([int]$row, 5).Value2 = $vals.E
}
